# issue #5: stock data output to json file
#
# Adds a "property_category" column (value "stock" for every data row) to
# the 股票 (stock) sheet, inserted right after the "total" column and
# before the "date" column. This pushes date / legislator_name /
# legislator_id one column to the right (H->I, I->J, J->K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H (column 8), shifting the existing H:J
# (date, legislator_name, legislator_id) columns right to I:K.
$ws.Columns.Item(8).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 8).Value = "property_category"

# Fill every existing data row (2 through 13) with the new category value.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
